$wb = $excel.ActiveWorkbook

# Sheet 1: ALC
$ws = $wb.Worksheets.Item(1)
# Row 43
$ws.Range("H43").Value = 1404095.6
$ws.Range("I43").Value = 1929656.8
$ws.Range("J43").Value = 2599.3333
$ws.Range("K43").Value = 1929656.8
$ws.Range("L43").Value = 2599.3333
$ws.Range("M43").Value = -1929587.8
$ws.Range("N43").Value = -2737.3333

# Row 111
$ws.Range("H111").Value = 1297.4
$ws.Range("I111").Value = 1297.4
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 3892.2
$ws.Range("L111").Value = 0
$ws.Range("M111").ClearContents()
$ws.Range("N111").Value = -825.2000000000003

# Row 132
$ws.Range("H132").Value = 7410.59
$ws.Range("I132").Value = 1956.196
$ws.Range("J132").Value = 13087.612
$ws.Range("K132").Value = 5868.588
$ws.Range("L132").Value = 39262.836
$ws.Range("M132").Value = -3338.588
$ws.Range("N132").Value = -44322.836

# Row 137
$ws.Range("H137").Value = 2411.6667
$ws.Range("I137").Value = 2690.1538
$ws.Range("J137").Value = 2230.65
$ws.Range("K137").Value = 8070.4614
$ws.Range("L137").Value = 6691.950000000001
$ws.Range("M137").Value = -5520.4614
$ws.Range("N137").Value = -11791.95

# Row 139
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").ClearContents()
$ws.Range("N139").Value = 0

# Row 141
$ws.Range("H141").Value = 6074.737
$ws.Range("I141").Value = 6201.1763
$ws.Range("K141").Value = 18603.5289
$ws.Range("M141").Value = -13423.5289

# Sheet 2: ARM
$ws = $wb.Worksheets.Item(2)
# Row 2
$ws.Range("H2").Value = 12226158
$ws.Range("I2").Value = 998005.9399999999
$ws.Range("K2").Value = 998005.9399999999
$ws.Range("M2").Value = -997892.9399999999

# Row 32
$ws.Range("H32").Value = 2968.5527
$ws.Range("I32").Value = 1617.6865
$ws.Range("J32").Value = 13025
$ws.Range("K32").Value = 1617.6865
$ws.Range("L32").Value = 13025
$ws.Range("M32").Value = -1330.6865
$ws.Range("N32").Value = -13599

# Row 63
$ws.Range("H63").Value = 4818
$ws.Range("I63").Value = 4780.1665
$ws.Range("J63").Value = 4874.75
$ws.Range("K63").Value = 4780.1665
$ws.Range("L63").Value = 4874.75
$ws.Range("M63").Value = -4094.1665
$ws.Range("N63").Value = -6246.75

# Row 66
$ws.Range("H66").Value = 4818
$ws.Range("I66").Value = 4780.1665
$ws.Range("J66").Value = 4874.75
$ws.Range("K66").Value = 23900.8325
$ws.Range("L66").Value = 24373.75
$ws.Range("M66").Value = -20468.8325
$ws.Range("N66").Value = -31237.75

# Row 74
$ws.Range("H74").Value = 7814717.5
$ws.Range("I74").Value = 11905755
$ws.Range("J74").Value = 4555.636
$ws.Range("K74").Value = 11905755
$ws.Range("L74").Value = 4555.636
$ws.Range("M74").Value = -11904881
$ws.Range("N74").Value = -6303.636

# Row 77
$ws.Range("H77").Value = 7814717.5
$ws.Range("I77").Value = 11905755
$ws.Range("J77").Value = 4555.636
$ws.Range("K77").Value = 59528775
$ws.Range("L77").Value = 22778.18
$ws.Range("M77").Value = -59524407
$ws.Range("N77").Value = -31514.18

# Row 102
$ws.Range("H102").Value = 1957620.1
$ws.Range("I102").Value = 1957620.1
$ws.Range("K102").Value = 1957620.1
$ws.Range("M102").Value = -1955998.1

# Row 116
$ws.Range("H116").Value = 12226158
$ws.Range("I116").Value = 998005.9399999999
$ws.Range("K116").Value = 998005.9399999999
$ws.Range("M116").Value = -995711.9399999999

# Row 132
$ws.Range("H132").Value = 50946.31
$ws.Range("I132").Value = 64698.668
$ws.Range("J132").Value = 20003.5
$ws.Range("K132").Value = 194096.004
$ws.Range("L132").Value = 60010.5
$ws.Range("M132").Value = -191566.004
$ws.Range("N132").Value = -65070.5

# Row 140
$ws.Range("H140").Value = 112714.5
$ws.Range("J140").Value = 112714.5
$ws.Range("L140").Value = 112714.5
$ws.Range("N140").Value = -123074.5

# Sheet 3: BSM
$ws = $wb.Worksheets.Item(3)
# Row 3
$ws.Range("H3").Value = 12226158
$ws.Range("I3").Value = 998005.9399999999
$ws.Range("K3").Value = 998005.9399999999
$ws.Range("M3").Value = -997891.9399999999

# Row 20
$ws.Range("H20").Value = 1782.1428
$ws.Range("I20").Value = 1586.9474
$ws.Range("K20").Value = 1586.9474
$ws.Range("M20").Value = -1339.9474

# Row 92
$ws.Range("H92").Value = 44499.25
$ws.Range("J92").Value = 44499.25
$ws.Range("L92").Value = 44499.25
$ws.Range("N92").Value = -49491.25

# Row 134
$ws.Range("H134").Value = 4870.567
$ws.Range("I134").Value = 3556
$ws.Range("J134").Value = 6589.615
$ws.Range("K134").Value = 10668
$ws.Range("L134").Value = 19768.845
$ws.Range("M134").Value = -8133
$ws.Range("N134").Value = -24838.845

# Row 140
$ws.Range("H140").Value = 162193
$ws.Range("J140").Value = 162193
$ws.Range("L140").Value = 162193
$ws.Range("N140").Value = -172553

# Sheet 4: CRP
$ws = $wb.Worksheets.Item(4)
# Row 109
$ws.Range("H109").Value = 56930
$ws.Range("J109").Value = 56930
$ws.Range("L109").Value = 56930
$ws.Range("N109").Value = -59010

# Row 115
$ws.Range("H115").Value = 53000
$ws.Range("J115").Value = 53000
$ws.Range("L115").Value = 53000
$ws.Range("N115").Value = -55350

# Row 141
$ws.Range("H141").Value = 211876
$ws.Range("J141").Value = 211876
$ws.Range("L141").Value = 211876
$ws.Range("N141").Value = -222236

# Sheet 5: CUL
$ws = $wb.Worksheets.Item(5)
# Row 68
$ws.Range("H68").Value = 174874.22
$ws.Range("J68").Value = 196491.34
$ws.Range("L68").Value = 589474.02
$ws.Range("N68").Value = -591096.02

# Row 71
$ws.Range("H71").Value = 174874.22
$ws.Range("J71").Value = 196491.34
$ws.Range("L71").Value = 1768422.06
$ws.Range("N71").Value = -1776534.06

# Row 98
$ws.Range("H98").Value = 1697.6842
$ws.Range("I98").Value = 1412.5
$ws.Range("J98").Value = 1773.7333
$ws.Range("K98").Value = 4237.5
$ws.Range("L98").Value = 5321.199900000001
$ws.Range("M98").Value = -2739.5
$ws.Range("N98").Value = -8317.1999

# Row 107
$ws.Range("H107").Value = 3695.8206
$ws.Range("J107").Value = 3867.879
$ws.Range("L107").Value = 11603.637
$ws.Range("N107").Value = -15443.637

# Row 129
$ws.Range("H129").Value = 3426.818
$ws.Range("I129").Value = 2000
$ws.Range("J129").Value = 3961.875
$ws.Range("K129").Value = 6000
$ws.Range("L129").Value = 11885.625
$ws.Range("M129").Value = -1000
$ws.Range("N129").Value = -21885.625

# Row 130
$ws.Range("H130").Value = 43750
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 43750
$ws.Range("K130").Value = 0
$ws.Range("L130").ClearContents()
$ws.Range("M130").Value = 131250
$ws.Range("N130").Value = -141290

# Row 131
$ws.Range("H131").Value = 18521288
$ws.Range("J131").Value = 6805717
$ws.Range("L131").Value = 20417151
$ws.Range("N131").Value = -20427231

# Sheet 6: GSM
$ws = $wb.Worksheets.Item(6)
# Row 70
$ws.Range("H70").Value = 749278.0600000001
$ws.Range("I70").Value = 1039960.75
$ws.Range("J70").Value = 6422.222
$ws.Range("K70").Value = 1039960.75
$ws.Range("L70").Value = 6422.222
$ws.Range("M70").Value = -1039690.75
$ws.Range("N70").Value = -6962.222

# Row 73
$ws.Range("H73").Value = 749278.0600000001
$ws.Range("I73").Value = 1039960.75
$ws.Range("J73").Value = 6422.222
$ws.Range("K73").Value = 1039960.75
$ws.Range("L73").Value = 6422.222
$ws.Range("M73").Value = -1039024.75
$ws.Range("N73").Value = -8294.222

# Row 102
$ws.Range("H102").Value = 378423.2
$ws.Range("I102").Value = 724118.9399999999
$ws.Range("K102").Value = 724118.9399999999
$ws.Range("M102").Value = -722496.9399999999

# Row 122
$ws.Range("H122").Value = 531099.5
$ws.Range("I122").Value = 1006099.6
$ws.Range("K122").Value = 3018298.8
$ws.Range("M122").Value = -3015848.8

# Row 132
$ws.Range("H132").Value = 6806.8125
$ws.Range("I132").Value = 6593.933
$ws.Range("K132").Value = 19781.799
$ws.Range("M132").Value = -17251.799

# Sheet 7: LTW
$ws = $wb.Worksheets.Item(7)
# Row 6
$ws.Range("H6").Value = 59498.25
$ws.Range("J6").Value = 59498.25
$ws.Range("L6").Value = 59498.25
$ws.Range("N6").Value = -59722.25

# Row 117
$ws.Range("H117").Value = 55000
$ws.Range("J117").Value = 55000
$ws.Range("L117").Value = 55000
$ws.Range("N117").Value = -64178

# Row 122
$ws.Range("H122").Value = 166682830
$ws.Range("I122").Value = 250003740
$ws.Range("J122").Value = 40975
$ws.Range("K122").Value = 750011220
$ws.Range("L122").Value = 122925
$ws.Range("M122").Value = -750008770
$ws.Range("N122").Value = -127825

# Row 123
$ws.Range("H123").Value = 58900
$ws.Range("J123").Value = 58900
$ws.Range("L123").Value = 58900
$ws.Range("N123").Value = -68700

# Row 129
$ws.Range("H129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("L129").ClearContents()
$ws.Range("N129").Value = 0

# Row 136
$ws.Range("H136").Value = 4034.05
$ws.Range("I136").Value = 4060.2444
$ws.Range("J136").Value = 3798.3
$ws.Range("K136").Value = 12180.7332
$ws.Range("L136").Value = 11394.9
$ws.Range("M136").Value = -9630.733200000001
$ws.Range("N136").Value = -16494.9

# Row 139
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").ClearContents()
$ws.Range("N139").Value = 0

# Sheet 8: WVR
$ws = $wb.Worksheets.Item(8)
# Row 107
$ws.Range("H107").Value = 1259.7273
$ws.Range("I107").Value = 1317.1111
$ws.Range("K107").Value = 3951.3333
$ws.Range("M107").Value = -2031.3333

# Row 115
$ws.Range("H115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("L115").ClearContents()
$ws.Range("N115").Value = 0

# Row 122
$ws.Range("H122").Value = 3136.5945
$ws.Range("I122").Value = 3187.3713
$ws.Range("K122").Value = 9562.1139
$ws.Range("M122").Value = -7112.1139

# Row 129
$ws.Range("H129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("L129").ClearContents()
$ws.Range("N129").Value = 0

# Row 132
$ws.Range("H132").Value = 20886246
$ws.Range("I132").Value = 59995.668
$ws.Range("J132").Value = 166670000
$ws.Range("K132").Value = 179987.004
$ws.Range("L132").Value = 500010000
$ws.Range("M132").Value = -177457.004
$ws.Range("N132").Value = -500015060

# Row 136
$ws.Range("H136").Value = 8531.790000000001
$ws.Range("I136").Value = 3752.3684
$ws.Range("J136").Value = 9652.888999999999
$ws.Range("K136").Value = 11257.1052
$ws.Range("L136").Value = 28958.667
$ws.Range("M136").Value = -8707.1052
$ws.Range("N136").Value = -34058.667

# Row 138
$ws.Range("H138").Value = 50000
$ws.Range("J138").Value = 50000
$ws.Range("L138").Value = 50000
$ws.Range("N138").Value = -60280
